$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "classes" sheet: fill in row 13 with the new "NoCardinalityClass" entry
#    (DEV-3755: support classes that have no cardinalities / property sheet)
# ---------------------------------------------------------------------------
$classes = $wb.Worksheets.Item("classes")
$classes.Range("A13").Value = "NoCardinalityClass"
$classes.Range("B13").Value = "Class Without Cardinalities"
$classes.Range("G13").Value = "No cardinalitiy"
$classes.Range("L13").Value = "Resource"
[void]$classes.Range("G13").Select()

# ---------------------------------------------------------------------------
# 2) "Owner" sheet: a blank row was inserted after "hasBirthDate" (row 12),
#    pushing the remaining property rows down by one, and the trailing
#    left-over "invalid" placeholder rows were cleaned up.
# ---------------------------------------------------------------------------
$owner = $wb.Worksheets.Item("Owner")
$owner.Rows.Item(13).Insert()
$owner.Range("A13").Value = "    "
$owner.Rows.Item(17).Delete()
$owner.Rows.Item(17).Delete()
$owner.Rows.Item(17).Delete()
[void]$owner.Range("A33").Select()

# ---------------------------------------------------------------------------
# 3) "GenericAnthroponym" sheet: cursor moved (no data change)
# ---------------------------------------------------------------------------
$generic = $wb.Worksheets.Item("GenericAnthroponym")
[void]$generic.Range("H44").Select()

[void]$classes.Activate()
